$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp (2024-08-12 -> 2024-08-13)
$ws.Name = "IClientBalance-20240813-095047-"

# Shift the reference date (column G, "Dt. Referencia") for every data row
# from 2024-08-12 (serial 45516) to 2024-08-13 (serial 45517)
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45517
}

# Updated balances for the rows whose figures changed in this refresh
# (column D = Vl. Projetado, E = Saldo Previsto, H = Vl. Total = D + E)
# Row 5
$ws.Cells.Item(5, 5).Value = 979.23
$ws.Cells.Item(5, 8).Value = 979.23
# Row 55
$ws.Cells.Item(55, 5).Value = 999.99
$ws.Cells.Item(55, 8).Value = 999.99
# Row 58
$ws.Cells.Item(58, 4).Value = 5179.03
$ws.Cells.Item(58, 8).Value = 5179.94
# Row 60
$ws.Cells.Item(60, 5).Value = 995.21
$ws.Cells.Item(60, 8).Value = 995.21
# Row 197
$ws.Cells.Item(197, 5).Value = 67.56
$ws.Cells.Item(197, 8).Value = 67.56
# Row 230
$ws.Cells.Item(230, 5).Value = 999.99
$ws.Cells.Item(230, 8).Value = 999.99
# Row 235
$ws.Cells.Item(235, 5).Value = 697.14
$ws.Cells.Item(235, 8).Value = 697.14
# Row 270
$ws.Cells.Item(270, 4).Value = -11960.75
$ws.Cells.Item(270, 5).Value = 12903.93
$ws.Cells.Item(270, 8).Value = 943.18
# Row 271
$ws.Cells.Item(271, 5).Value = 986.93
$ws.Cells.Item(271, 8).Value = 986.93
